# Apply "Trade #65 closed" update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.3    # Current Capital
$summary.Range("B4").Value = 0.3       # Total P&L $
$summary.Range("B5").Value = 0.09      # Total P&L %
$summary.Range("B6").Value = 65        # Total Trades
$summary.Range("B8").Value = 34        # Losing Trades
$summary.Range("B9").Value = 32.31     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.3      # Capital
$status.Range("D4").Value = 65         # Trades
$status.Range("E4").Value = 0.3        # P&L $
$status.Range("F4").Value = 0.3        # P&L %
$status.Range("G4").Value = 32.31      # Win Rate %

# ---------------------------------------------------------------------------
# New trade row data (Trade #65), appended to both "All Trades" and
# "MarketMaking" sheets as row 66.
# ---------------------------------------------------------------------------
$newRow = @(65, "2026-02-17", "15:47:26", "MarketMaking", "UP", 0.65, 0.57, "CLOSED", -12.3077, -0.08, 100.3, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.12)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Columns B (Date) and C (Time) are stored as plain text in the source
    # data, not native dates/times. Force text format first so Excel does
    # not auto-convert "2026-02-17" / "15:47:26" into date/time serials.
    $ws.Range("B66:C66").NumberFormat = "@"

    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item(66, $i + 1).Value = $newRow[$i]
    }
}
